# "Fix minor typo in Ex 3 slide deck"
#
# 1) Slide 12 ("Text Placeholder 2"): the run break in
#    "Build your First  ListView" is wrong (a stray double-space before
#    "ListView", and the word "First" split off onto the wrong run). Re-split
#    the text into "Build your " + "First ListView" so the rendered text
#    reads "Build your First ListView" instead of "Build your First  ListView".
#
# 2) The cached date/time field text ("3/3/2019 9:45 AM" -> "12/1/2019 5:24 PM")
#    on the handout master, notes master and every notes page is also
#    refreshed to match. (Those placeholders hold a live <a:fld type="datetime8">
#    field — PowerPoint recalculates/caches its display text automatically;
#    we still try to nudge it here via the documented object model and quietly
#    continue if a given host/master/notes shape does not allow a scripted
#    rewrite of that cached field text.)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Part 1: slide 12 text fix
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(12)
$shape = $slide.Shapes.Item(3)
$tr = $shape.TextFrame.TextRange

$fullText = $tr.Text
$oldRun1 = "Build your First  "
$oldRun2 = "ListView"
$newRun1 = "Build your "
$newRun2 = "First ListView"

$hit = $fullText.IndexOf($oldRun1 + $oldRun2)
if ($hit -ge 0) {
    $start = $hit + 1   # TextRange.Characters is 1-based

    # Update the second run first so the first run's character offsets
    # (which we already know) stay valid while we still need them.
    $run2 = $tr.Characters($start + $oldRun1.Length, $oldRun2.Length)
    $run2.Text = $newRun2

    $run1 = $tr.Characters($start, $oldRun1.Length)
    $run1.Text = $newRun1
}

# ---------------------------------------------------------------------------
# Part 2: refresh the cached "3/3/2019 9:45 AM" -> "12/1/2019 5:24 PM"
# datetime8 field text wherever it appears (handout master, notes master,
# and each slide's notes page). Best-effort / non-fatal.
# ---------------------------------------------------------------------------
$oldDate = "3/3/2019 9:45 AM"
$newDate = "12/1/2019 5:24 PM"

function Set-DateShapeText($shapeCollection) {
    for ($i = 1; $i -le $shapeCollection.Count; $i++) {
        $shp = $shapeCollection.Item($i)
        try {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                $txt = $shp.TextFrame.TextRange.Text
                if ($txt -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        } catch {
            # This host does not allow rewriting this placeholder's text
            # (e.g. master/notes date fields) - leave the auto field as-is.
        }
    }
}

try {
    Set-DateShapeText $p.HandoutMaster.Shapes
} catch { }

try {
    if ($p.HasNotesMaster) {
        Set-DateShapeText $p.NotesMaster.Shapes
    }
} catch { }

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    try {
        $sld = $p.Slides.Item($si)
        if ($sld.HasNotesPage) {
            Set-DateShapeText $sld.NotesPage.Shapes
        }
    } catch { }
}
